$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the four new "architecture diagram" paragraphs right after the
#    "High-level software architecture description..." heading and before
#    the paragraph that holds the embedded Visio drawing (OLE object).
# ---------------------------------------------------------------------------

# locate the heading paragraph reliably via the Paragraphs collection
# (Range.Text carries a trailing CR for the paragraph mark, strip it first)
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($txt -eq "High-level software architecture description, including the communication interfaces between the different elements") {
        $target = $i
        break
    }
}

$headingParagraph = $d.Paragraphs($target)
$insertionPoint = $d.Range($headingParagraph.Range.End, $headingParagraph.Range.End)

$para1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">The following architecture diagram </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>gives an overview of the proposed solution and can be divided in</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>to</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> four main parts. First, there are the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>IoT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Devices, whose general functionality have been described in the previous chapter. There will be some additional information about each device in the following chapters. Some </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>IoT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Devices don&#8217;t directly push the obtained data to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>IoT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Agent, so a special reader/collector is needed. This can be seen in the architecture containing several RFID </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Ea</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>rtags</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> and RFID Readers as well.</w:t></w:r></w:p>
'@

$para2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Furthermore, there exist several different </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>IoT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Agents, that manage the communication with the mentioned </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>IoT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Devices.</w:t></w:r></w:p>
'@

$para3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">All the collected data will be handed over to the Context Broker, which is the central part of the architecture. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>In order to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> store this data for statistical and long-term analysis </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>of the obtained</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> data, Cygnus and a Data Analyser in Hadoop is used.</w:t></w:r></w:p>
'@

$para4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Finally, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">there are three different types of User Applications, that get the information of interest from the Context Broker. The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Wirecloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Webpage helps the farmer to get an overview of all the different components and data to be processed. With the smartphone and computer applications the farmer can easily access whatever information he needs (current information of certain cows, long-term information, statistics etc.) and the farmer will be notified and alerted in case of emergency. In case of emergency the Stable </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Notifier</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> helps to quickly detect the affected cattle. </w:t></w:r></w:p>
'@

$trailerEmpty = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>'

$fullXml = $para1 + $para2 + $para3 + $para4 + $trailerEmpty
$insertionPoint.InsertXML($fullXml)

# the insertion leaves one stray empty paragraph right before the OLE
# paragraph (because the "last" xml fragment always gets merged into
# whatever paragraph used to sit at the insertion point) -- remove it.
for ($i = $target + 1; $i -le $target + 6; $i++) {
    $p = $d.Paragraphs($i)
    $ptxt = $p.Range.Text.TrimEnd([char]13)
    if ($ptxt -eq "") {
        $p.Range.Delete()
        break
    }
}

Write-Output "Inserted architecture paragraphs."
